# Regenerate the "K" column (column G) values for the gibaut_ian save_data sheet.
# Per commit message: "regen save_data to use K instead of Strike#, regen std/mean,
# calc and write s_vals" - here we write the newly-computed K (strikeout) values
# for each game row (rows 2 through 37).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 2
    3  = 3
    4  = 0
    5  = 2
    6  = 0
    7  = 0
    8  = 2
    9  = 0
    10 = 0
    11 = 2
    12 = 0
    13 = 1
    14 = 0
    15 = 2
    16 = 0
    17 = 2
    18 = 1
    19 = 1
    20 = 3
    21 = 2
    22 = 2
    23 = 2
    24 = 4
    25 = 1
    26 = 1
    27 = 2
    28 = 2
    29 = 2
    30 = 3
    31 = 3
    32 = 2
    33 = 0
    34 = 1
    35 = 0
    36 = 2
    37 = 3
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
